# "render termino de version entregable"
# Update the summary row (row 2) of the Resumen_global sheet with the
# final/closing values for pct_outliers-related metrics:
#   A2: 8.538 -> 8.535
#   B2: 9     -> 6
#   C2: 0,11  -> 0,07
#   G2: 13,29 -> 10,33
#
# These are text values (Spanish-locale formatted numbers using "," as the
# decimal separator and "." as the thousands separator) stored as shared
# strings, not numeric cells. Assigning the literal text directly via
# .Value would make Excel's locale-aware parser coerce strings such as
# "8.535" or "6" into real numbers, changing the cell type. To keep them
# as genuine text cells (matching the original file), each new value is
# entered as a literal-text formula (="...") and then converted in place
# to a plain value via Copy / Paste Special (values only), which leaves a
# plain text cell behind with no residual formula and no style changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Formula = '="8.535"'
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B2").Formula = '="6"'
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("C2").Formula = '="0,07"'
$ws.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("G2").Formula = '="10,33"'
$ws.Range("G2").Copy()
$ws.Range("G2").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = $false
